$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from column CG (the last existing data column) to column CH
# so the new column matches the look/style of its neighbour.
$ws.Range("CG1:CG11").Copy()
$ws.Range("CH1:CH11").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("CH1").Value = "8-nov"
$ws.Range("CH2").Value = 7
$ws.Range("CH3").Value = 11
$ws.Range("CH4").Value = 7
$ws.Range("CH5").Value = 11
$ws.Range("CH6").Value = 7
$ws.Range("CH7").Value = 8
$ws.Range("CH8").Value = 10
$ws.Range("CH9").Value = 11
$ws.Range("CH10").Value = 8
$ws.Range("CH11").Value = 0

$ws.Range("CM19").Select() | Out-Null
